# Insert two new data rows at row 362 (pushing the existing rows 362-425
# down to 364-427), then populate the two new rows with their data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("362:363").Insert()

# Row 362: Mandarina, Tango, Especial — Provincia de Quillota
$ws.Cells.Item(362, 1).Value = 9
$ws.Cells.Item(362, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(362, 3).Value = "Metropolitana"
$ws.Cells.Item(362, 4).Value = 44522
$ws.Cells.Item(362, 5).Value = 13
$ws.Cells.Item(362, 6).Value = "Fruta"
$ws.Cells.Item(362, 7).Value = 100102
$ws.Cells.Item(362, 8).Value = "Cítricos"
$ws.Cells.Item(362, 9).Value = 100102004
$ws.Cells.Item(362, 10).Value = "Mandarina"
$ws.Cells.Item(362, 11).Value = "Tango"
$ws.Cells.Item(362, 12).Value = "Especial"
$ws.Cells.Item(362, 13).Value = 220
$ws.Cells.Item(362, 14).Value = 10800
$ws.Cells.Item(362, 15).Value = 10800
$ws.Cells.Item(362, 16).Value = 10800
$ws.Cells.Item(362, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(362, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(362, 19).Value = 600
$ws.Cells.Item(362, 20).Value = 18

# Row 363: Mandarina, Tango, Primera — Provincia de Quillota
$ws.Cells.Item(363, 1).Value = 9
$ws.Cells.Item(363, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(363, 3).Value = "Metropolitana"
$ws.Cells.Item(363, 4).Value = 44522
$ws.Cells.Item(363, 5).Value = 13
$ws.Cells.Item(363, 6).Value = "Fruta"
$ws.Cells.Item(363, 7).Value = 100102
$ws.Cells.Item(363, 8).Value = "Cítricos"
$ws.Cells.Item(363, 9).Value = 100102004
$ws.Cells.Item(363, 10).Value = "Mandarina"
$ws.Cells.Item(363, 11).Value = "Tango"
$ws.Cells.Item(363, 12).Value = "Primera"
$ws.Cells.Item(363, 13).Value = 200
$ws.Cells.Item(363, 14).Value = 9000
$ws.Cells.Item(363, 15).Value = 9000
$ws.Cells.Item(363, 16).Value = 9000
$ws.Cells.Item(363, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(363, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(363, 19).Value = 500
$ws.Cells.Item(363, 20).Value = 18
